$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bug report date corrected (one day earlier)
$ws.Range("B2").Value = 45017

# Status changed from "Новый" (New) to "выполнено" (done)
$ws.Range("B9").Value = "выполнено"

# Severity changed from "critical" to "major"
$ws.Range("B10").Value = "major"

# "Включить ландшафтную ориентацию" step cell gets centered + wrapped alignment
$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("C13").VerticalAlignment = -4108
$ws.Range("C13").WrapText = $true

# Update the active selection to reflect the edited cell
$ws.Range("B8").Select()
